$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Programming 1" / "Programming 2-7" rows are being merged into a single
# "Programming" row. Delete the "Programming 1" row (old row 2); the
# "Programming 2-7" row (old row 3) slides up to row 2 and everything below
# it (Quiz, Discussion Board, Demo, Final Project, Totals) shifts up by one.
$ws.Rows("2:2").Delete() | Out-Null

# Rework the row that slid into position 2 into the combined "Programming" row.
$ws.Range("A2").Value = "Programming"
$ws.Range("B2").Value = 7
$ws.Range("C2").Value = 32

# "Quiz" (now row 3): update counts.
$ws.Range("B3").Value = 8
$ws.Range("C3").Value = 16

# "Discussion Board" (row 4) and "Demo" (row 5) are unchanged.

# "Final Project" (now row 6): total points changes from 32 to 36.
$ws.Range("C6").Value = 36

# Refresh the sheet's selection to the new used range (A1:E7).
$ws.Range("A1:E7").Select() | Out-Null
